$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
# Delete rows 5 and 6 (shifts all subsequent rows up by two).
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(5).Delete()

# --- Sheet 2: "Monthly Trend" ---
# Update B4 from 360 to 140.
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B4").Value = 140
